# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for the "Artic Snow" Nectarín variety
# (Femacal de La Calera market, week of 2022-03-08 / serial 44628) just
# above the existing historical block, pushing the old rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 889; Excel copies the
# formatting (incl. the date style on column D) from the row above.
$ws.Rows("889:892").Insert()

# Row 889: Artic Snow - Especial
$ws.Range("A889").Value = 3
$ws.Range("B889").Value = "Femacal de La Calera"
$ws.Range("C889").Value = "Coquimbo"
$ws.Range("D889").Value = 44628
$ws.Range("E889").Value = 5
$ws.Range("F889").Value = "Fruta"
$ws.Range("G889").Value = 100103
$ws.Range("H889").Value = "Frutos de hueso (carozo)"
$ws.Range("I889").Value = 100103006
$ws.Range("J889").Value = "Nectarín"
$ws.Range("K889").Value = "Artic Snow"
$ws.Range("L889").Value = "Especial"
$ws.Range("M889").Value = 78
$ws.Range("N889").Value = 15000
$ws.Range("O889").Value = 15000
$ws.Range("P889").Value = 15000
$ws.Range("Q889").Value = "$/caja 15 kilos empedrada"
$ws.Range("R889").Value = "Región de O'Higgins"
$ws.Range("S889").Value = 1000
$ws.Range("T889").Value = 15

# Row 890: Artic Snow - Extra (doble especial)
$ws.Range("A890").Value = 3
$ws.Range("B890").Value = "Femacal de La Calera"
$ws.Range("C890").Value = "Coquimbo"
$ws.Range("D890").Value = 44628
$ws.Range("E890").Value = 5
$ws.Range("F890").Value = "Fruta"
$ws.Range("G890").Value = 100103
$ws.Range("H890").Value = "Frutos de hueso (carozo)"
$ws.Range("I890").Value = 100103006
$ws.Range("J890").Value = "Nectarín"
$ws.Range("K890").Value = "Artic Snow"
$ws.Range("L890").Value = "Extra (doble especial)"
$ws.Range("M890").Value = 70
$ws.Range("N890").Value = 17000
$ws.Range("O890").Value = 17000
$ws.Range("P890").Value = 17000
$ws.Range("Q890").Value = "$/caja 15 kilos empedrada"
$ws.Range("R890").Value = "Región de O'Higgins"
$ws.Range("S890").Value = 1133
$ws.Range("T890").Value = 15

# Row 891: Artic Snow - Primera
$ws.Range("A891").Value = 3
$ws.Range("B891").Value = "Femacal de La Calera"
$ws.Range("C891").Value = "Coquimbo"
$ws.Range("D891").Value = 44628
$ws.Range("E891").Value = 5
$ws.Range("F891").Value = "Fruta"
$ws.Range("G891").Value = 100103
$ws.Range("H891").Value = "Frutos de hueso (carozo)"
$ws.Range("I891").Value = 100103006
$ws.Range("J891").Value = "Nectarín"
$ws.Range("K891").Value = "Artic Snow"
$ws.Range("L891").Value = "Primera"
$ws.Range("M891").Value = 70
$ws.Range("N891").Value = 13000
$ws.Range("O891").Value = 13000
$ws.Range("P891").Value = 13000
$ws.Range("Q891").Value = "$/caja 15 kilos empedrada"
$ws.Range("R891").Value = "Región de O'Higgins"
$ws.Range("S891").Value = 867
$ws.Range("T891").Value = 15

# Row 892: Artic Snow - Segunda
$ws.Range("A892").Value = 3
$ws.Range("B892").Value = "Femacal de La Calera"
$ws.Range("C892").Value = "Coquimbo"
$ws.Range("D892").Value = 44628
$ws.Range("E892").Value = 5
$ws.Range("F892").Value = "Fruta"
$ws.Range("G892").Value = 100103
$ws.Range("H892").Value = "Frutos de hueso (carozo)"
$ws.Range("I892").Value = 100103006
$ws.Range("J892").Value = "Nectarín"
$ws.Range("K892").Value = "Artic Snow"
$ws.Range("L892").Value = "Segunda"
$ws.Range("M892").Value = 70
$ws.Range("N892").Value = 11000
$ws.Range("O892").Value = 11000
$ws.Range("P892").Value = 11000
$ws.Range("Q892").Value = "$/caja 15 kilos empedrada"
$ws.Range("R892").Value = "Región de O'Higgins"
$ws.Range("S892").Value = 733
$ws.Range("T892").Value = 15
